$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark additional tasks as "Completed" (status column updates)
$ws.Range("D17").Value = "Completed"
$ws.Range("D18").Value = "Completed"
$ws.Range("D20").Value = "Completed"
$ws.Range("D25").Value = "Completed"
$ws.Range("D26").Value = "Completed"
$ws.Range("D29").Value = "Completed"

# Turn on the table's AutoFilter dropdowns
$lo = $ws.ListObjects.Item(1)
$lo.ShowAutoFilter = $true

# Move the active selection back to the top of the Status column
$ws.Range("D1").Select() | Out-Null
